$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 92 (hunk 0)
$ws.Range("H92").Value = 1573.55
$ws.Range("I92").Value = 1691.1428
$ws.Range("K92").Value = 1691.1428
$ws.Range("M92").Value = -443.1428000000001

# row 125 (hunk 1)
$ws.Range("H125").Value = 111112310
$ws.Range("I125").Value = 200000420
$ws.Range("K125").Value = 1800003780
$ws.Range("M125").Value = -1800001320

$ws = $wb.Worksheets.Item("ARM")
# row 32 (hunk 2)
$ws.Range("H32").Value = 2966.0205
$ws.Range("I32").Value = 2519.0344
$ws.Range("J32").Value = 6501.273
$ws.Range("K32").Value = 2519.0344
$ws.Range("L32").Value = 6501.273
$ws.Range("M32").Value = -2232.0344
$ws.Range("N32").Value = -7075.273

# row 45 (hunk 3)
$ws.Range("H45").Value = 1137.909
$ws.Range("I45").Value = 1147
$ws.Range("J45").Value = 1130.3334
$ws.Range("K45").Value = 1147
$ws.Range("L45").Value = 1130.3334
$ws.Range("M45").Value = -770
$ws.Range("N45").Value = -1884.3334

# row 74 (hunk 4)
$ws.Range("H74").Value = 865.2564
$ws.Range("I74").Value = 840.8333
$ws.Range("J74").Value = 946.6667
$ws.Range("K74").Value = 840.8333
$ws.Range("L74").Value = 946.6667
$ws.Range("M74").Value = 33.16669999999999
$ws.Range("N74").Value = -2694.6667

# row 77 (hunk 5)
$ws.Range("H77").Value = 865.2564
$ws.Range("I77").Value = 840.8333
$ws.Range("J77").Value = 946.6667
$ws.Range("K77").Value = 4204.1665
$ws.Range("L77").Value = 4733.3335
$ws.Range("M77").Value = 163.8334999999997
$ws.Range("N77").Value = -13469.3335

# row 88 (hunk 6)
$ws.Range("H88").Value = 1255915
$ws.Range("J88").Value = 3462.8
$ws.Range("L88").Value = 3462.8
$ws.Range("N88").Value = -4274.8

# row 91 (hunk 7)
$ws.Range("H91").Value = 1255915
$ws.Range("J91").Value = 3462.8
$ws.Range("L91").Value = 3462.8
$ws.Range("N91").Value = -6270.8

# row 132 (hunk 8)
$ws.Range("H132").Value = 3545.0715
$ws.Range("I132").Value = 3915.3096
$ws.Range("J132").Value = 2434.3572
$ws.Range("K132").Value = 11745.9288
$ws.Range("L132").Value = 7303.071599999999
$ws.Range("M132").Value = -9215.9288
$ws.Range("N132").Value = -12363.0716

$ws = $wb.Worksheets.Item("BSM")
# row 103 (hunk 9)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0

$ws = $wb.Worksheets.Item("CRP")
# row 5 (hunk 10)
$ws.Range("H5").Value = 337.5
$ws.Range("I5").Value = 253.4
$ws.Range("J5").Value = 421.6
$ws.Range("K5").Value = 253.4
$ws.Range("L5").Value = 421.6
$ws.Range("M5").Value = -141.4
$ws.Range("N5").Value = -645.6

# row 31 (hunk 11)
$ws.Range("H31").Value = 6948234.5
$ws.Range("I31").Value = 3132
$ws.Range("K31").Value = 3132
$ws.Range("M31").Value = -2837

# row 34 (hunk 12)
$ws.Range("H34").Value = 6948234.5
$ws.Range("I34").Value = 3132
$ws.Range("K34").Value = 3132
$ws.Range("M34").Value = -2930

# row 58 (hunk 13)
$ws.Range("H58").Value = 805.2759
$ws.Range("I58").Value = 679.7308
$ws.Range("J58").Value = 1893.3334
$ws.Range("K58").Value = 679.7308
$ws.Range("L58").Value = 1893.3334
$ws.Range("M58").Value = -476.7308
$ws.Range("N58").Value = -2299.3334

# row 94 (hunk 14)
$ws.Range("H94").Value = 47619480
$ws.Range("J94").Value = 515
$ws.Range("L94").Value = 515
$ws.Range("N94").Value = -1417

# row 99 (hunk 15)
$ws.Range("H99").Value = 36348.62
$ws.Range("I99").Value = 60539.06
$ws.Range("J99").Value = 2078.8333
$ws.Range("K99").Value = 60539.06
$ws.Range("L99").Value = 2078.8333
$ws.Range("M99").Value = -59041.06
$ws.Range("N99").Value = -5074.8333

# row 126 (hunk 16)
$ws.Range("H126").Value = 36348.62
$ws.Range("I126").Value = 60539.06
$ws.Range("J126").Value = 2078.8333
$ws.Range("K126").Value = 181617.18
$ws.Range("L126").Value = 6236.499899999999
$ws.Range("M126").Value = -179147.18
$ws.Range("N126").Value = -11176.4999

# row 132 (hunk 17)
$ws.Range("H132").Value = 3003
$ws.Range("I132").Value = 2360.3076
$ws.Range("K132").Value = 7080.9228
$ws.Range("M132").Value = -4550.9228

# row 134 (hunk 18)
$ws.Range("H134").Value = 899.61224
$ws.Range("I134").Value = 836.2105
$ws.Range("J134").Value = 1118.6364
$ws.Range("K134").Value = 2508.6315
$ws.Range("L134").Value = 3355.9092
$ws.Range("M134").Value = 26.36850000000004
$ws.Range("N134").Value = -8425.9092

# row 136 (hunk 19)
$ws.Range("H136").Value = 805.2759
$ws.Range("I136").Value = 679.7308
$ws.Range("J136").Value = 1893.3334
$ws.Range("K136").Value = 2039.1924
$ws.Range("L136").Value = 5680.0002
$ws.Range("M136").Value = 510.8075999999999
$ws.Range("N136").Value = -10780.0002

$ws = $wb.Worksheets.Item("CUL")
# row 5 (hunk 20)
$ws.Range("H5").Value = 666.56525
$ws.Range("I5").Value = 391.85715
$ws.Range("J5").Value = 1093.8889
$ws.Range("K5").Value = 1175.57145
$ws.Range("L5").Value = 3281.6667
$ws.Range("M5").Value = -1063.57145
$ws.Range("N5").Value = -3505.6667

# row 131 (hunk 21)
$ws.Range("H131").Value = 2648803
$ws.Range("J131").Value = 4630513
$ws.Range("L131").Value = 13891539
$ws.Range("N131").Value = -13901619

# row 135 (hunk 22)
$ws.Range("H135").Value = 666.56525
$ws.Range("I135").Value = 391.85715
$ws.Range("J135").Value = 1093.8889
$ws.Range("K135").Value = 3526.71435
$ws.Range("L135").Value = 9845.000099999999
$ws.Range("M135").Value = -991.7143499999997
$ws.Range("N135").Value = -14915.0001

$ws = $wb.Worksheets.Item("GSM")
# row 102 (hunk 23)
$ws.Range("H102").Value = 1214.625
$ws.Range("I102").Value = 1110.3077
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 1110.3077
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = 511.6922999999999
$ws.Range("N102").Value = -4910.6666

# row 126 (hunk 24)
$ws.Range("H126").Value = 2410.182
$ws.Range("I126").Value = 5420.6665
$ws.Range("J126").Value = 1281.25
$ws.Range("K126").Value = 16261.9995
$ws.Range("L126").Value = 3843.75
$ws.Range("M126").Value = -13791.9995
$ws.Range("N126").Value = -8783.75

# row 132 (hunk 25)
$ws.Range("H132").Value = 36897.93
$ws.Range("I132").Value = 45264.26
$ws.Range("K132").Value = 135792.78
$ws.Range("M132").Value = -133262.78

$ws = $wb.Worksheets.Item("LTW")
# row 40 (hunk 26)
$ws.Range("H40").Value = 2003.9
$ws.Range("I40").Value = 1654.875
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 1654.875
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -1518.875
$ws.Range("N40").Value = -3672

# row 132 (hunk 27)
$ws.Range("H132").Value = 2929.889
$ws.Range("I132").Value = 3061.111
$ws.Range("K132").Value = 9183.332999999999
$ws.Range("M132").Value = -6653.332999999999

# row 136 (hunk 28)
$ws.Range("H136").Value = 5662.0415
$ws.Range("I136").Value = 9970.362999999999
$ws.Range("J136").Value = 2016.5385
$ws.Range("K136").Value = 29911.089
$ws.Range("L136").Value = 6049.6155
$ws.Range("M136").Value = -27361.089
$ws.Range("N136").Value = -11149.6155

$ws = $wb.Worksheets.Item("WVR")
# row 107 (hunk 29)
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1800
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5640

# row 126 (hunk 30)
$ws.Range("H126").Value = 7029.1577
$ws.Range("I126").Value = 8812.214
$ws.Range("J126").Value = 2036.6
$ws.Range("K126").Value = 26436.642
$ws.Range("L126").Value = 6109.799999999999
$ws.Range("M126").Value = -23966.642
$ws.Range("N126").Value = -11049.8

